$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ibbs-web-app-metrics")
$ws2 = $wb.Worksheets.Item("Summary")

[void]$ws1.Range("A2:J19").Clear()

[void]$ws2.Activate()
[void]$ws2.Range("AT3").Select()

[void]$ws1.Activate()
[void]$ws1.Range("A2").Select()

Write-Host "done"
